$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 80
$ws1.Range("F4").Value = 241
$ws1.Range("F6").Value = 9988
$ws1.Range("F7").Value = 904
$ws1.Range("F9").Value = 1243
$ws1.Range("F10").Value = 5021
$ws1.Range("F11").Value = 8
$ws1.Range("F12").Value = 9
$ws1.Range("F13").Value = 181
$ws1.Range("F15").Value = 68
$ws1.Range("F17").Value = 295
$ws1.Range("F18").Value = 587
$ws1.Range("F21").Value = 15
$ws1.Range("F22").Value = 1507

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 80
$ws4.Range("F5").Value = 241
$ws4.Range("F7").Value = 9988
$ws4.Range("F8").Value = 904
$ws4.Range("F10").Value = 1243
$ws4.Range("F11").Value = 5021
$ws4.Range("F12").Value = 8
$ws4.Range("F13").Value = 9
$ws4.Range("F14").Value = 181
$ws4.Range("F15").Value = 120
$ws4.Range("F16").Value = 68
$ws4.Range("F18").Value = 295
$ws4.Range("F19").Value = 587
$ws4.Range("F22").Value = 15
$ws4.Range("F23").Value = 1507
